$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("M21").Select()
